$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6588
$ws.Range("F3").Value = 768
$ws.Range("F4").Value = 1102
$ws.Range("F5").Value = 115
$ws.Range("F6").Value = 630
$ws.Range("F7").Value = 214
$ws.Range("F8").Value = 45
$ws.Range("F9").Value = 807
$ws.Range("F10").Value = 1266
$ws.Range("F11").Value = 28
$ws.Range("F12").Value = 100
$ws.Range("F13").Value = 517
$ws.Range("F14").Value = 507
$ws.Range("F15").Value = 357
$ws.Range("F17").Value = 1455
$ws.Range("F19").Value = 429
$ws.Range("F22").Value = 1095
$ws.Range("F23").Value = 211
$ws.Range("F24").Value = 2290
$ws.Range("F26").Value = 157
$ws.Range("F29").Value = 3725
$ws.Range("F31").Value = 685

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 34
$ws.Range("F11").Value = 133
$ws.Range("F17").Value = 388
$ws.Range("F19").Value = 4107
$ws.Range("F21").Value = 15
$ws.Range("F25").Value = 236
$ws.Range("F29").Value = 38

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1225
$ws.Range("F5").Value = 1606
$ws.Range("F8").Value = 911

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1225
$ws.Range("F4").Value = 1606
$ws.Range("F7").Value = 911
$ws.Range("F8").Value = 6588
$ws.Range("F9").Value = 34
$ws.Range("F10").Value = 768
$ws.Range("F12").Value = 115
$ws.Range("F13").Value = 630
$ws.Range("F14").Value = 214
$ws.Range("F15").Value = 45
$ws.Range("F16").Value = 807
$ws.Range("F19").Value = 133
$ws.Range("F20").Value = 133
$ws.Range("F23").Value = 1266
$ws.Range("F24").Value = 28
$ws.Range("F25").Value = 100
$ws.Range("F26").Value = 517
$ws.Range("F27").Value = 507
$ws.Range("F29").Value = 388
$ws.Range("F31").Value = 357
$ws.Range("F33").Value = 1455
$ws.Range("F36").Value = 429
$ws.Range("F40").Value = 1095
$ws.Range("F41").Value = 211
$ws.Range("F42").Value = 2290
$ws.Range("F45").Value = 157
$ws.Range("F47").Value = 3725
$ws.Range("F51").Value = 685
